$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.357.32'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').Value = '2.547.77'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('D5').Value = '''592.58'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '''175.73'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''0.525'
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('D9').Value = '2.549.71'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').Value = '  -2.38%  '
$ws.Range('D11').Value = '''0.167'
$ws.Range('E11').Value = '  +1.80%  '
$ws.Range('D12').Value = '''0.345'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('E13').Value = '  -2.38%  '
$ws.Range('D14').Value = '''26.58'
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').Value = '2.982.68'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = '68.407.57'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.536.68'
$ws.Range('E18').Value = '  +1.08%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''11.96'
$ws.Range('E19').Value = '  +4.46%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''8.01'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('B21').Value = 'Binance-PegBSC-USD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D21').Value = '''1.70'
$ws.Range('E21').Value = '  +69.93%  '
$ws.Range('D22').Value = '''366.91'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('D23').Value = '''4.17'
$ws.Range('E23').Value = '  -0.55%  '
$ws.Range('D24').Value = '''4.57'
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('D25').Value = '''72.11'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = '''1.90'
$ws.Range('E27').Value = '  -4.16%  '
$ws.Range('D28').Value = '''9.95'
$ws.Range('E28').Value = '  -3.37%  '
$ws.Range('D29').Value = '2.681.06'
$ws.Range('E29').Value = '  +1.01%  '
$ws.Range('D30').Value = '0.0₃0966'
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D31').Value = '''535.89'
$ws.Range('E31').Value = '  -3.39%  '
$ws.Range('D32').Value = '''8.33'
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').Value = '''1.31'
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('E34').Value = '  +0.59%  '
$ws.Range('D35').Value = '''0.129'
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').Value = '''159.89'
$ws.Range('E37').Value = '  +2.53%  '
$ws.Range('E38').Value = '  -2.10%  '
$ws.Range('E39').Value = '  +3.10%  '
$ws.Range('D40').Value = '''18.64'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '''1.79'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').Value = '''5.13'
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('D43').Value = '''0.348'
$ws.Range('E43').Value = '  -2.07%  '
$ws.Range('D44').Value = '''2.51'
$ws.Range('E44').Value = '  -0.61%  '
$ws.Range('D45').Value = '''0.997'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').Value = '''39.42'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('D47').Value = '''149.00'
$ws.Range('E47').Value = '  +1.13%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '''0.556'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').Value = '''3.71'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0278'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('E51').Value = '  +1.73%  '
